# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" worksheet right before "总计" (i.e. right after
#    "2021-Q4") with the per-fund holdings for that quarter.
# 2) Insert a new leading row into the "总计" (summary) sheet with the
#    2022-Q1 totals, pushing the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper cells used purely as formatting donors via Copy/PasteSpecial:
#   $blank  -> a cell that has never been touched, i.e. carries the
#              workbook's default (no "s" attribute) style.
#   $styled -> a cell that already carries the bold-font + boxed-border
#              "header / index column" style used on the "2021-Q4" and
#              "总计" sheets.
# ---------------------------------------------------------------------
$blank = $wb.Worksheets.Item("2021-Q3").Range("A1")
$styled = $wb.Worksheets.Item("2021-Q4").Range("B1")

# Write $text into $range as a genuine text value (never auto-coerced to
# a number, so things like fund codes keep their leading zeros and
# "9.55" stays text rather than becoming the number 9.55), ending up with
# the workbook's plain default style (no borders/bold).
function Set-Text($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $blank.Copy()
    $range.PasteSpecial(-4122)
}

# Same as Set-Text but also stamps the bold/bordered "header" style onto
# the range afterwards (re-using the existing style rather than creating
# a new one).
function Set-StyledText($range, [string]$text) {
    Set-Text $range $text
    $styled.Copy()
    $range.PasteSpecial(-4122)
}

# Plain numeric value, default style.
function Set-Number($range, $number) {
    $range.Value = $number
}

# Plain numeric value with the bold/bordered "index column" style.
function Set-StyledNumber($range, $number) {
    $styled.Copy()
    $range.PasteSpecial(-4122)
    $range.Value = $number
}

# ===========================================================================
# 1) New "2022-Q1" worksheet
# ===========================================================================
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $afterSheet)
$q1.Name = "2022-Q1"

# Match the outlinePr (summaryBelow/summaryRight) settings used by every
# other sheet in this workbook.
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1

# ---- Header row ----
Set-StyledText $q1.Range("B1") "基金代码"
Set-StyledText $q1.Range("C1") "基金名称"
Set-StyledText $q1.Range("D1") "基金规模"
Set-StyledText $q1.Range("E1") "股票总仓位"
Set-StyledText $q1.Range("F1") "仓位占比"
Set-StyledText $q1.Range("G1") "持有市值(亿元)"
Set-StyledText $q1.Range("H1") "仓位排名"

# ---- Data rows ----
Set-StyledNumber $q1.Range("A2") 0
Set-Text $q1.Range("B2") "003298"
Set-Text $q1.Range("C2") "嘉实物流产业股票A"
Set-Text $q1.Range("D2") "9.55"
Set-Text $q1.Range("E2") "86.21"
Set-Text $q1.Range("F2") "2.92"
Set-Text $q1.Range("G2") "0.2789"
Set-Number $q1.Range("H2") 9

Set-StyledNumber $q1.Range("A3") 1
Set-Text $q1.Range("B3") "003299"
Set-Text $q1.Range("C3") "嘉实物流产业股票C"
Set-Text $q1.Range("D3") "4.23"
Set-Text $q1.Range("E3") "86.21"
Set-Text $q1.Range("F3") "2.92"
Set-Text $q1.Range("G3") "0.1235"
Set-Number $q1.Range("H3") 9

Set-StyledNumber $q1.Range("A4") 2
Set-Text $q1.Range("B4") "005459"
Set-Text $q1.Range("C4") "银河嘉谊灵活配置混合A"
Set-Text $q1.Range("D4") "6.47"
Set-Text $q1.Range("E4") "39.69"
Set-Text $q1.Range("F4") "0.53"
Set-Text $q1.Range("G4") "0.0343"
Set-Number $q1.Range("H4") 10

Set-StyledNumber $q1.Range("A5") 3
Set-Text $q1.Range("B5") "005460"
Set-Text $q1.Range("C5") "银河嘉谊灵活配置混合C"
Set-Text $q1.Range("D5") "2.79"
Set-Text $q1.Range("E5") "39.69"
Set-Text $q1.Range("F5") "0.53"
Set-Text $q1.Range("G5") "0.0148"
Set-Number $q1.Range("H5") 10

Set-StyledNumber $q1.Range("A6") 4
Set-Text $q1.Range("B6") "004250"
Set-Text $q1.Range("C6") "银河量化优选混合"
Set-Text $q1.Range("D6") "0.39"
Set-Text $q1.Range("E6") "80.03"
Set-Text $q1.Range("F6") "1.67"
Set-Text $q1.Range("G6") "0.0065"
Set-Number $q1.Range("H6") 7

# ===========================================================================
# 2) "总计" summary sheet: insert a new leading data row for 2022-Q1
# ===========================================================================
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Clean up whatever formatting Insert() copied down onto the new row so we
# can apply exactly the styles we want below.
$blank.Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

Set-StyledNumber $total.Range("A2") 0
Set-Text $total.Range("B2") "2022-Q1"
Set-Number $total.Range("C2") 5
Set-Number $total.Range("D2") 0.46

# Insert() shifted the old rows 2/3 down to 3/4 verbatim, so the leading
# "index" column (A) still holds its old 0 / 1 values - renumber it to
# stay a sequential 0-based index for the new row order.
Set-Number $total.Range("A3") 1
Set-Number $total.Range("A4") 2

# ===========================================================================
# Restore the originally active sheet / clear clipboard marching ants.
# ===========================================================================
$wb.Worksheets.Item("2021-Q3").Activate()
$excel.CutCopyMode = $false
